# Question option 1 mandatory removed
#
# - E1 ("QuestionOption1" header) loses its "mandatory" highlight fill; it
#   should match the plain formatting already used by the other optional
#   question-option headers (F1:I1).
# - A new "Question Type" value ("Question with description-WITHD") is
#   filled in for every data row (R2:R9) in the already-present "Question
#   Type" column R (its data-validation dropdown already existed).
# - The view's scroll position / active selection are updated to reflect
#   where the editor ended up after making the change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Question")

# Remove the mandatory-field highlight from E1 by copying F1's (unhighlighted)
# formatting onto it - this reuses the existing "normal" cell style instead of
# fabricating a new one, so E1 ends up styled identically to F1:I1.
$ws.Range("F1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill the new "Question Type" column (R) for every data row with the
# "Question with description-WITHD" option.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 18).Value = "Question with description-WITHD"
}

# Reflect the editor's final scroll position and selection.
$win = $excel.Windows.Item(1)
$win.ScrollRow = 3
$win.ScrollColumn = 10
$ws.Range("R13").Select()
